# Daily attendance processing - reorders the "Recorded By" (column G) list
# so the most recently recording user/system appears first.
#
# For every row, if the current value of column G exactly matches one of the
# known "before" combinations, replace it with the corresponding reordered
# combination. Rows whose value isn't one of these combinations (e.g. a
# single recorder, or "System, admin@admin.com") are left untouched, exactly
# matching the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "System, dnasr281@gmail.com" = "dnasr281@gmail.com, System";
    "System, backup@backdoor.com" = "backup@backdoor.com, System";
    "System, system, backup@backdoor.com" = "backup@backdoor.com, System, system";
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com";
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Range("G" + $row)
    $val = $cell.Value()
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
